$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.458.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.506.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.79%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.68'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.09'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.79%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.504.29'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.80%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.64'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.969.86'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.55%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.402.33'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.526.65'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.15%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.23'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.51'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.64'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.01'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.03'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0978'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '528.66'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.50%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.87%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.45'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.74'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.27%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.52'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.07%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.57'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.36'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.95%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.51%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0272'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -9.11%  '
